$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B149").Value = 48654
$ws.Range("E149").Value = 38.26
$ws.Range("F149").Value = -1
$ws.Range("G149").Value = -32.02
$ws.Range("B150").Value = 63902
$ws.Range("E150").Value = 34.04
$ws.Range("F150").Value = 2
$ws.Range("G150").Value = 64.04000000000001
$ws.Range("B161").Value = 57756
$ws.Range("F161").Value = -100
$ws.Range("G161").Value = -6644
$ws.Range("B162").Value = 53925
$ws.Range("E162").Value = 79.37
$ws.Range("F162").Value = 1
$ws.Range("G162").Value = 66.44
$ws.Range("B163").Value = 64350
$ws.Range("E163").Value = 70.63
$ws.Range("F163").Value = 101
$ws.Range("G163").Value = 6710.44
$ws.Range("B183").Value = 57552
$ws.Range("E183").Value = 136.86
$ws.Range("F183").Value = -5
$ws.Range("G183").Value = -603.45
$ws.Range("B184").Value = 64329
$ws.Range("E184").Value = 128.32
$ws.Range("F184").Value = 6
$ws.Range("G184").Value = 724.14
$ws.Range("B264").Value = 48719
$ws.Range("E264").Value = 353.35
$ws.Range("F264").Value = -81
$ws.Range("G264").Value = -23955.75
$ws.Range("B265").Value = 64979
$ws.Range("E265").Value = 314.41
$ws.Range("F265").Value = 82
$ws.Range("G265").Value = 24251.5
$ws.Range("B317").Value = 61610
$ws.Range("E317").Value = 122.71
$ws.Range("F317").Value = -58
$ws.Range("G317").Value = -5957.18
$ws.Range("B318").Value = 63565
$ws.Range("E318").Value = 109.19
$ws.Range("F318").Value = 60
$ws.Range("G318").Value = 6162.6
$ws.Range("B346").Value = 55373
$ws.Range("E346").Value = 163.62
$ws.Range("F346").Value = -94
$ws.Range("G346").Value = -13562.32
$ws.Range("B347").Value = 63520
$ws.Range("E347").Value = 153.4
$ws.Range("F347").Value = 97
$ws.Range("G347").Value = 13995.16
$ws.Range("B350").Value = 57802
$ws.Range("E350").Value = 162.71
$ws.Range("F350").Value = -79
$ws.Range("G350").Value = -11334.92
$ws.Range("B352").Value = 63571
$ws.Range("E352").Value = 152.53
$ws.Range("F352").Value = 27
$ws.Range("G352").Value = 3873.96
$ws.Range("B355").Value = 63510
$ws.Range("E355").Value = 50.66
$ws.Range("F355").Value = 167
$ws.Range("G355").Value = 7955.88
$ws.Range("B356").Value = 55356
$ws.Range("E356").Value = 54.04
$ws.Range("F356").Value = -158
$ws.Range("G356").Value = -7527.12
$ws.Range("B375").Value = 61605
$ws.Range("E375").Value = 133.78
$ws.Range("F375").Value = -13
$ws.Range("G375").Value = -1455.48
$ws.Range("B376").Value = 63563
$ws.Range("E376").Value = 119.04
$ws.Range("F376").Value = 15
$ws.Range("G376").Value = 1679.4
$ws.Range("B382").Value = 60325
$ws.Range("E382").Value = 151.57
$ws.Range("F382").Value = -102
$ws.Range("G382").Value = -12939.72
$ws.Range("B383").Value = 63560
$ws.Range("E383").Value = 134.87
$ws.Range("F383").Value = 104
$ws.Range("G383").Value = 13193.44
$ws.Range("B389").Value = 57817
$ws.Range("F389").Value = 3
$ws.Range("G389").Value = 239.43
$ws.Range("B390").Value = 62865
$ws.Range("F390").Value = 151
$ws.Range("G390").Value = 12051.31
$ws.Range("B400").Value = 62933
$ws.Range("F400").Value = 146
$ws.Range("G400").Value = 8632.98
$ws.Range("B401").Value = 57835
$ws.Range("F401").Value = 1
$ws.Range("G401").Value = 59.13
$ws.Range("B419").Value = 57856
$ws.Range("F419").Value = 2
$ws.Range("G419").Value = 342.66
$ws.Range("B420").Value = 63007
$ws.Range("F420").Value = 984
$ws.Range("G420").Value = 168588.72
$ws.Range("B431").Value = 63102
$ws.Range("C431").Value = 'HUL-Vim Bar Multipack Fw 4X200G'
$ws.Range("F431").Value = 36
$ws.Range("G431").Value = 2140.92
$ws.Range("B432").Value = 53082
$ws.Range("C432").Value = 'HUL-VIM BAR MULTIPACK FW 4X200G'
$ws.Range("F432").Value = 1
$ws.Range("G432").Value = 59.47
$ws.Range("B586").Value = 45695
$ws.Range("E586").Value = 23.58
$ws.Range("F586").Value = -36
$ws.Range("G586").Value = -710.28
$ws.Range("B587").Value = 64915
$ws.Range("E587").Value = 20.98
$ws.Range("F587").Value = 40
$ws.Range("G587").Value = 789.2
$ws.Range("B590").Value = 64922
$ws.Range("E590").Value = 20.98
$ws.Range("F590").Value = 207
$ws.Range("G590").Value = 4084.11
$ws.Range("B591").Value = 45706
$ws.Range("E591").Value = 23.58
$ws.Range("F591").Value = -202
$ws.Range("G591").Value = -3985.46
$ws.Range("B593").Value = 45718
$ws.Range("E593").Value = 19.38
$ws.Range("F593").Value = -294
$ws.Range("G593").Value = -4768.68
$ws.Range("B594").Value = 64927
$ws.Range("E594").Value = 17.26
$ws.Range("F594").Value = 295
$ws.Range("G594").Value = 4784.9
$ws.Range("B687").Value = 53319
$ws.Range("E687").Value = 310.64
$ws.Range("F687").Value = -6
$ws.Range("G687").Value = -1643.52
$ws.Range("B688").Value = 64810
$ws.Range("E688").Value = 291.22
$ws.Range("F688").Value = 7
$ws.Range("G688").Value = 1917.44
